$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for each data row.
# All data rows (2 to 28) need their date bumped by one day: 45174 -> 45175.
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
